$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the original inlineStr type) instead of
# auto-converting to a numeric value.
$numericLookingCells = @("D5","D6","D10","D12","D14","D19","D21","D23","D24","D27","D28","D30","D31","D32","D33","D35","D41","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value changes from the diff
$ws.Range('D2').Value = '44.473.55'
$ws.Range('E2').Value = '  +3.67%  '
$ws.Range('D3').Value = '2.290.67'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '320.04'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').Value = '106.07'
$ws.Range('E6').Value = '  +7.90%  '
$ws.Range('E7').Value = '  +1.83%  '
$ws.Range('E9').Value = '  +2.71%  '
$ws.Range('D10').Value = '39.37'
$ws.Range('E10').Value = '  +7.71%  '
$ws.Range('E11').Value = '  +2.45%  '
$ws.Range('D12').Value = '7.97'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('E13').Value = '  +2.18%  '
$ws.Range('D14').Value = '0.888'
$ws.Range('E14').Value = '  +2.95%  '
$ws.Range('D15').Value = '2.635.73'
$ws.Range('E15').Value = '  +3.32%  '
$ws.Range('E16').Value = '  +4.03%  '
$ws.Range('D17').Value = '2.295.82'
$ws.Range('E17').Value = '  +4.00%  '
$ws.Range('D18').Value = '44.362.55'
$ws.Range('E18').Value = '  +3.91%  '
$ws.Range('D19').Value = '14.28'
$ws.Range('E19').Value = '  -3.23%  '
$ws.Range('E20').Value = '  +4.53%  '
$ws.Range('D21').Value = '6.60'
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('E22').Value = '  +2.16%  '
$ws.Range('D23').Value = '3.23'
$ws.Range('E23').Value = '  +2.66%  '
$ws.Range('D24').Value = '238.92'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('E25').Value = '  +4.66%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '10.36'
$ws.Range('E27').Value = '  +2.85%  '
$ws.Range('D28').Value = '39.46'
$ws.Range('E28').Value = '  +16.34%  '
$ws.Range('D30').Value = '6.61'
$ws.Range('E30').Value = '  +5.48%  '
$ws.Range('D31').Value = '164.12'
$ws.Range('E31').Value = '  +5.40%  '
$ws.Range('D32').Value = '0.0891'
$ws.Range('E32').Value = '  +2.35%  '
$ws.Range('D33').Value = '20.62'
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('D35').Value = '3.33'
$ws.Range('E35').Value = '  +4.61%  '
$ws.Range('E36').Value = '  +5.49%  '
$ws.Range('E37').Value = '  +13.38%  '
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('E39').Value = '  +2.78%  '
$ws.Range('E40').Value = '  +8.03%  '
$ws.Range('D41').Value = '15.64'
$ws.Range('E41').Value = '  +29.30%  '
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').Value = '1.775.69'
$ws.Range('E44').Value = '  -5.24%  '
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('D46').Value = '86.19'
$ws.Range('E46').Value = '  -3.66%  '
$ws.Range('D47').Value = '5.43'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '76.54'
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '8.94'
$ws.Range('E49').Value = '  +3.95%  '
$ws.Range('D50').Value = '60.04'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('D51').Value = '105.27'
$ws.Range('E51').Value = '  +4.13%  '

# Restore default (Normal) style on the cells we temporarily reformatted,
# so no stray number-format styling is left behind on them.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
